$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 690.5238000000001
$ws.Range("I33").Value = 559.125
$ws.Range("K33").Value = 559.125
$ws.Range("M33").Value = -330.125
# Row 40
$ws.Range("H40").Value = 5565850.5
$ws.Range("J40").Value = 12363691
$ws.Range("L40").Value = 12363691
$ws.Range("N40").Value = -12364041
# Row 58
$ws.Range("H58").Value = 22876
$ws.Range("I58").Value = 612
$ws.Range("J58").Value = 100800
$ws.Range("K58").Value = 1836
$ws.Range("L58").Value = 302400
$ws.Range("M58").Value = -1686
$ws.Range("N58").Value = -302700
# Row 97
$ws.Range("H97").Value = 3147.125
$ws.Range("J97").Value = 3147.125
$ws.Range("L97").Value = 9441.375
$ws.Range("N97").Value = -10433.375
# Row 98
$ws.Range("H98").Value = 3366.6428
$ws.Range("I98").Value = 1558.2727
$ws.Range("K98").Value = 1558.2727
$ws.Range("M98").Value = -60.27269999999999
# Row 122
$ws.Range("H122").Value = 3366.6428
$ws.Range("I122").Value = 1558.2727
$ws.Range("K122").Value = 4674.8181
$ws.Range("M122").Value = -2224.8181
# Row 125
$ws.Range("H125").Value = 10928203
$ws.Range("I125").Value = 4239034
$ws.Range("K125").Value = 38151306
$ws.Range("M125").Value = -38148846
# Row 126
$ws.Range("H126").Value = 92497
$ws.Range("J126").Value = 92497
$ws.Range("L126").Value = 92497
$ws.Range("N126").Value = -102377
# Row 127
$ws.Range("H127").Value = 737.6
$ws.Range("I127").Value = 737.6
$ws.Range("K127").Value = 2212.8
$ws.Range("M127").Value = 2747.2
# Row 132
$ws.Range("H132").Value = 3436.3
$ws.Range("I132").Value = 3436.3
$ws.Range("K132").Value = 10308.9
$ws.Range("M132").Value = -7778.900000000001
# Row 137
$ws.Range("H137").Value = 3008.2
$ws.Range("I137").Value = 2564.6667
$ws.Range("K137").Value = 7694.000100000001
$ws.Range("M137").Value = -5144.000100000001

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 52633320
$ws.Range("I61").Value = 58825264
$ws.Range("K61").Value = 58825264
$ws.Range("M61").Value = -58825052
# Row 74
$ws.Range("H74").Value = 31254510
$ws.Range("I74").Value = 33338092
$ws.Range("K74").Value = 33338092
$ws.Range("M74").Value = -33337218
# Row 77
$ws.Range("H77").Value = 31254510
$ws.Range("I77").Value = 33338092
$ws.Range("K77").Value = 166690460
$ws.Range("M77").Value = -166686092
# Row 132
$ws.Range("H132").Value = 3848862.2
$ws.Range("I132").Value = 4350236
$ws.Range("J132").Value = 4998
$ws.Range("K132").Value = 13050708
$ws.Range("L132").Value = 14994
$ws.Range("M132").Value = -13048178
$ws.Range("N132").Value = -20054
# Row 136
$ws.Range("H136").Value = 52633320
$ws.Range("I136").Value = 58825264
$ws.Range("K136").Value = 176475792
$ws.Range("M136").Value = -176473242

$ws = $wb.Worksheets.Item("BSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").Value = ""
# Row 88
$ws.Range("H88").Value = 61499
$ws.Range("I88").Value = 57998
$ws.Range("K88").Value = 57998
$ws.Range("M88").Value = -57592
# Row 91
$ws.Range("H91").Value = 61499
$ws.Range("I91").Value = 57998
$ws.Range("K91").Value = 57998
$ws.Range("M91").Value = -56594
# Row 107
$ws.Range("H107").Value = 147201.42
$ws.Range("I107").Value = 5068.5
$ws.Range("K107").Value = 5068.5
$ws.Range("M107").Value = -3148.5
# Row 134
$ws.Range("H134").Value = 45460204
$ws.Range("I134").Value = 50006096
$ws.Range("K134").Value = 150018288
$ws.Range("M134").Value = -150015753

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 398
$ws.Range("I7").Value = 547
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 547
$ws.Range("L7").Value = 100
$ws.Range("M7").Value = -434
$ws.Range("N7").Value = -326
# Row 99
$ws.Range("H99").Value = 3139.9285
$ws.Range("I99").Value = 3042.4167
$ws.Range("K99").Value = 3042.4167
$ws.Range("M99").Value = -1544.4167
# Row 126
$ws.Range("H126").Value = 3139.9285
$ws.Range("I126").Value = 3042.4167
$ws.Range("K126").Value = 9127.250100000001
$ws.Range("M126").Value = -6657.250100000001
# Row 132
$ws.Range("I132").Value = 333339680
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 1000019040
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -1000016510
$ws.Range("N132").Value = ""
# Row 134
$ws.Range("H134").Value = 9658200
$ws.Range("I134").Value = 9658200
$ws.Range("K134").Value = 28974600
$ws.Range("M134").Value = -28972065

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 407.875
$ws.Range("I2").Value = 239.28572
$ws.Range("J2").Value = 643.9
$ws.Range("K2").Value = 1435.71432
$ws.Range("L2").Value = 3863.4
$ws.Range("M2").Value = -1322.71432
$ws.Range("N2").Value = -4089.4
# Row 4
$ws.Range("H4").Value = 196811.83
$ws.Range("I4").Value = 217791.38
$ws.Range("K4").Value = 653374.14
$ws.Range("M4").Value = -653262.14

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 93.9375
$ws.Range("I2").Value = 26.333334
$ws.Range("J2").Value = 180.85715
$ws.Range("K2").Value = 26.333334
$ws.Range("L2").Value = 180.85715
$ws.Range("M2").Value = 86.66666599999999
$ws.Range("N2").Value = -406.85715
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
# Row 113
$ws.Range("H113").Value = 66199.19
$ws.Range("I113").Value = 86949.336
$ws.Range("J113").Value = 3948.75
$ws.Range("K113").Value = 86949.336
$ws.Range("L113").Value = 3948.75
$ws.Range("M113").Value = -84779.336
$ws.Range("N113").Value = -8288.75
# Row 122
$ws.Range("H122").Value = 107084.5
$ws.Range("J122").Value = 12003.2
$ws.Range("L122").Value = 36009.60000000001
$ws.Range("N122").Value = -40909.60000000001
# Row 132
$ws.Range("H132").Value = 7356383
$ws.Range("I132").Value = 7815970.5
$ws.Range("J132").Value = 2980
$ws.Range("K132").Value = 23447911.5
$ws.Range("L132").Value = 8940
$ws.Range("M132").Value = -23445381.5
$ws.Range("N132").Value = -14000

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4045.625
$ws.Range("I7").Value = 4045.625
$ws.Range("K7").Value = 4045.625
$ws.Range("M7").Value = -3933.625
# Row 22
$ws.Range("H22").Value = 3049.9285
$ws.Range("I22").Value = 3381.7273
$ws.Range("J22").Value = 1833.3334
$ws.Range("K22").Value = 3381.7273
$ws.Range("L22").Value = 1833.3334
$ws.Range("M22").Value = -3086.7273
$ws.Range("N22").Value = -2423.3334
# Row 27
$ws.Range("H27").Value = 3049.9285
$ws.Range("I27").Value = 3381.7273
$ws.Range("J27").Value = 1833.3334
$ws.Range("K27").Value = 3381.7273
$ws.Range("L27").Value = 1833.3334
$ws.Range("M27").Value = -3274.7273
$ws.Range("N27").Value = -2047.3334
# Row 93
$ws.Range("H93").Value = 2951
$ws.Range("I93").Value = 2951
$ws.Range("K93").Value = 2951
$ws.Range("M93").Value = -1703
# Row 126
$ws.Range("H126").Value = 4045.625
$ws.Range("I126").Value = 4045.625
$ws.Range("K126").Value = 12136.875
$ws.Range("M126").Value = -9666.875
# Row 132
$ws.Range("H132").Value = 15631763
$ws.Range("I132").Value = 19236800
$ws.Range("J132").Value = 9931.333000000001
$ws.Range("K132").Value = 57710400
$ws.Range("L132").Value = 29793.999
$ws.Range("M132").Value = -57707870
$ws.Range("N132").Value = -34853.999

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 2011.238
$ws.Range("I126").Value = 2133.923
$ws.Range("K126").Value = 6401.768999999999
$ws.Range("M126").Value = -3931.768999999999
